$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Existing data goes from row 2 to row 115 (A = value-2, B = 1).
# We need to extend it down to row 202 (A values 114..200, B = 1).

$startRow = 116
$endRow = 202

# Copy the style/format of the last existing data row (115) down to the new rows
$srcRange = $ws.Range("A115:B115")
$destRange = $ws.Range("A116:B202")
$srcRange.Copy()
$destRange.PasteSpecial(-4122)  # xlPasteFormats

for ($r = $startRow; $r -le $endRow; $r++) {
    $aValue = $r - 2
    $ws.Cells.Item($r, 1).Value = $aValue
    $ws.Cells.Item($r, 2).Value = 1
}

$excel.CutCopyMode = 0
